$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns from right to left so indices don't shift unexpectedly:
# Remove K (Email Mahasiswa), I (Long Title), G (Enrichment Track), E (Program), B (Binusian ID), A (Campus)
$ws.Range("K1").EntireColumn.Delete()
$ws.Range("I1").EntireColumn.Delete()
$ws.Range("G1").EntireColumn.Delete()
$ws.Range("E1").EntireColumn.Delete()
$ws.Range("B1").EntireColumn.Delete()
$ws.Range("A1").EntireColumn.Delete()

$ws.Range("N11").Select()
